$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the shared string "ID_transacao" -> "ID_Transacao" (header in A1)
$ws.Range("A1").Value = "ID_Transacao"

# 2. Shift column A values down by one (1..5 -> 0..4) for existing rows 2-6
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

# 3. Prepare formatting for the two new rows (7 and 8) by copying the
#    formatting of the existing last data row (row 6), which already has
#    the correct centered/date styles applied.
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A6:F6").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("G6").Copy()
$ws.Range("G7").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("G6").Copy()
$ws.Range("G8").PasteSpecial(-4122) # xlPasteFormats

# 4. Add new row 7: A7=6, B7=1, C7=1, D7=1, E7=1, F7=0, G7=45964 (2025-11-03)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 45964

# 5. Add new row 8: A8=5, B8=1, C8=1, D8=1, E8=1, F8=0, G8=45964 (2025-11-03)
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 45964

# 6. Update selection to C10 as in the diff
$ws.Range("C10").Select()
